$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 2072
$ws.Range('K3').Value = 8179
$ws.Range('L3').Value = 2087
$ws.Range('L4').Value = 579
$ws.Range('L5').Value = 120
$ws.Range('L6').Value = 1873
$ws.Range('K7').Value = 27553
$ws.Range('L7').Value = 6731

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L3').Value = 18
$ws.Range('L7').Value = 82

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L3').Value = 146
$ws.Range('L5').Value = 18
$ws.Range('L7').Value = 427

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 76
$ws.Range('L6').Value = 103
$ws.Range('L7').Value = 299

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L2').Value = 39
$ws.Range('L7').Value = 100

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 77
$ws.Range('L3').Value = 68
$ws.Range('L7').Value = 240

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 49
$ws.Range('L7').Value = 130

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L3').Value = 45
$ws.Range('L7').Value = 104

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L7').Value = 223
$ws.Range('L8').Value = 427
$ws.Range('L17').Value = 12
$ws.Range('L20').Value = 175
$ws.Range('L23').Value = 67
$ws.Range('L25').Value = 36
$ws.Range('L29').Value = 349
$ws.Range('L33').Value = 299
$ws.Range('L36').Value = 96
$ws.Range('L37').Value = 240
$ws.Range('L40').Value = 16
$ws.Range('L41').Value = 32
$ws.Range('L42').Value = 210
$ws.Range('L44').Value = 46
$ws.Range('L47').Value = 47
$ws.Range('K49').Value = 154
$ws.Range('L52').Value = 137
$ws.Range('L53').Value = 82
$ws.Range('L54').Value = 141
$ws.Range('L55').Value = 61
$ws.Range('L56').Value = 3
$ws.Range('K63').Value = 88
$ws.Range('L65').Value = 130
$ws.Range('L67').Value = 245
$ws.Range('K70').Value = 51
$ws.Range('L73').Value = 54
$ws.Range('L75').Value = 27
$ws.Range('L85').Value = 357
$ws.Range('L86').Value = 53
$ws.Range('L88').Value = 93
$ws.Range('K89').Value = 411
$ws.Range('L89').Value = 87
$ws.Range('L91').Value = 94
$ws.Range('L94').Value = 81
$ws.Range('L95').Value = 100
$ws.Range('L96').Value = 62
$ws.Range('L97').Value = 61
$ws.Range('L98').Value = 50
$ws.Range('L99').Value = 104
$ws.Range('K101').Value = 27553
$ws.Range('L101').Value = 6731

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 70
$ws.Range('L3').Value = 81
$ws.Range('L4').Value = 22
$ws.Range('L7').Value = 245

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K3').Value = 34
$ws.Range('K7').Value = 154

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L3').Value = 28
$ws.Range('L7').Value = 141

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 112
$ws.Range('L3').Value = 125
$ws.Range('L6').Value = 93
$ws.Range('L7').Value = 349

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('L2').Value = 18
$ws.Range('L6').Value = 14
$ws.Range('L7').Value = 46

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('L4').Value = 3
$ws.Range('L7').Value = 32

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L3').Value = 58
$ws.Range('L7').Value = 210

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L3').Value = 22
$ws.Range('L6').Value = 27

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L2').Value = 25
$ws.Range('L7').Value = 61

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L5').Value = 2
$ws.Range('L7').Value = 67

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L2').Value = 27
$ws.Range('L6').Value = 15
$ws.Range('L7').Value = 62

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L3').Value = 31
$ws.Range('L7').Value = 94

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L3').Value = 53
$ws.Range('L4').Value = 14
$ws.Range('L7').Value = 175

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('L3').Value = 4
$ws.Range('L7').Value = 12

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L4').Value = 7
$ws.Range('L7').Value = 96

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 65
$ws.Range('L7').Value = 223

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L3').Value = 21
$ws.Range('L4').Value = 11
$ws.Range('L7').Value = 81

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('L2').Value = 11
$ws.Range('L7').Value = 36

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L2').Value = 18
$ws.Range('L4').Value = 4
$ws.Range('L7').Value = 47

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('L6').Value = 27
$ws.Range('L7').Value = 50

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L2').Value = 19
$ws.Range('L6').Value = 14
$ws.Range('L7').Value = 54

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('L6').Value = 37
$ws.Range('L7').Value = 61

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('K4').Value = 9
$ws.Range('K7').Value = 51

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('L2').Value = 24
$ws.Range('L3').Value = 31
$ws.Range('L7').Value = 93

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K4').Value = 49
$ws.Range('L6').Value = 20
$ws.Range('K7').Value = 411
$ws.Range('L7').Value = 87

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L4').Value = 31
$ws.Range('L7').Value = 53

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('L2').Value = 16
$ws.Range('L7').Value = 27

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 110
$ws.Range('L6').Value = 65
$ws.Range('L7').Value = 357

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range('L6').Value = 2
$ws.Range('L7').Value = 3

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('L3').Value = 8
$ws.Range('L7').Value = 16

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L3').Value = 42
$ws.Range('L6').Value = 38
$ws.Range('L7').Value = 137
